$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "36.446.42"
$ws.Cells.Item(2, 5).Value = "  -2.80%  "

Set-TextValue $ws.Cells.Item(3, 4) "1.983.05"
$ws.Cells.Item(3, 5).Value = "  -3.52%  "

$ws.Cells.Item(4, 5).Value = "  +0.09%  "

Set-TextValue $ws.Cells.Item(5, 4) "244.54"
$ws.Cells.Item(5, 5).Value = "  -3.24%  "

Set-TextValue $ws.Cells.Item(6, 4) "0.628"
$ws.Cells.Item(6, 5).Value = "  -4.06%  "

Set-TextValue $ws.Cells.Item(7, 4) "59.45"
$ws.Cells.Item(7, 5).Value = "  -9.14%  "

$ws.Cells.Item(8, 5).Value = "  +0.08%  "

Set-TextValue $ws.Cells.Item(9, 4) "0.375"
$ws.Cells.Item(9, 5).Value = "  -2.39%  "

Set-TextValue $ws.Cells.Item(10, 4) "57.66"
$ws.Cells.Item(10, 5).Value = "  -3.61%  "

Set-TextValue $ws.Cells.Item(11, 4) "0.0820"
$ws.Cells.Item(11, 5).Value = "  +6.49%  "

$ws.Cells.Item(12, 5).Value = "  -1.09%  "

Set-TextValue $ws.Cells.Item(13, 4) "23.83"
$ws.Cells.Item(13, 5).Value = "  +4.40%  "

Set-TextValue $ws.Cells.Item(14, 4) "0.865"
$ws.Cells.Item(14, 5).Value = "  -6.05%  "

Set-TextValue $ws.Cells.Item(15, 4) "14.01"
$ws.Cells.Item(15, 5).Value = "  -5.92%  "

Set-TextValue $ws.Cells.Item(16, 4) "2.273.33"
$ws.Cells.Item(16, 5).Value = "  -3.48%  "

Set-TextValue $ws.Cells.Item(17, 4) "5.47"
$ws.Cells.Item(17, 5).Value = "  -2.18%  "

Set-TextValue $ws.Cells.Item(18, 4) "1.983.47"
$ws.Cells.Item(18, 5).Value = "  -3.59%  "

Set-TextValue $ws.Cells.Item(19, 4) "36.388.56"
$ws.Cells.Item(19, 5).Value = "  -2.60%  "

Set-TextValue $ws.Cells.Item(20, 4) "70.16"
$ws.Cells.Item(20, 5).Value = "  -4.80%  "

Set-TextValue $ws.Cells.Item(21, 4) "0.0₃0863"
$ws.Cells.Item(21, 5).Value = "  -1.76%  "

Set-TextValue $ws.Cells.Item(22, 4) "5.32"
$ws.Cells.Item(22, 5).Value = "  -3.10%  "

Set-TextValue $ws.Cells.Item(23, 4) "234.18"
$ws.Cells.Item(23, 5).Value = "  -2.51%  "

$ws.Cells.Item(24, 5).Value = "  +0.05%  "

$ws.Cells.Item(25, 5).Value = "  -1.69%  "

$ws.Cells.Item(26, 5).Value = "  -4.08%  "

Set-TextValue $ws.Cells.Item(27, 4) "10.03"
$ws.Cells.Item(27, 5).Value = "  -1.35%  "

Set-TextValue $ws.Cells.Item(28, 4) "162.41"

Set-TextValue $ws.Cells.Item(29, 4) "19.81"
$ws.Cells.Item(29, 5).Value = "  -1.12%  "

$ws.Cells.Item(30, 5).Value = "  +10.72%  "

$ws.Cells.Item(31, 5).Value = "  -2.35%  "

$ws.Cells.Item(32, 5).Value = "  -1.35%  "

$ws.Cells.Item(33, 5).Value = "  -7.07%  "

Set-TextValue $ws.Cells.Item(34, 4) "0.0630"
$ws.Cells.Item(34, 5).Value = "  +0.46%  "

$ws.Cells.Item(35, 5).Value = "  -6.13%  "

Set-TextValue $ws.Cells.Item(36, 4) "6.31"
$ws.Cells.Item(36, 5).Value = "  +4.29%  "

$ws.Cells.Item(37, 5).Value = "  +0.05%  "

Set-TextValue $ws.Cells.Item(38, 4) "2.26"
$ws.Cells.Item(38, 5).Value = "  -7.67%  "

$ws.Cells.Item(39, 5).Value = "  -2.70%  "

$ws.Cells.Item(40, 5).Value = "  +1.69%  "

$ws.Cells.Item(41, 5).Value = "  -0.89%  "

Set-TextValue $ws.Cells.Item(42, 4) "0.0963"
$ws.Cells.Item(42, 5).Value = "  -7.45%  "

$ws.Cells.Item(43, 5).Value = "  -3.86%  "

$ws.Cells.Item(44, 5).Value = "  -2.34%  "

$ws.Cells.Item(45, 5).Value = "  -5.17%  "

Set-TextValue $ws.Cells.Item(46, 4) "92.68"
$ws.Cells.Item(46, 5).Value = "  -4.39%  "

Set-TextValue $ws.Cells.Item(47, 4) "16.22"
$ws.Cells.Item(47, 5).Value = "  -6.08%  "

Set-TextValue $ws.Cells.Item(48, 4) "1.375.55"
$ws.Cells.Item(48, 5).Value = "  -3.44%  "

Set-TextValue $ws.Cells.Item(49, 4) "7.50"
$ws.Cells.Item(49, 5).Value = "  -6.01%  "

$ws.Cells.Item(50, 5).Value = "  -3.29%  "

Set-TextValue $ws.Cells.Item(51, 4) "45.30"
$ws.Cells.Item(51, 5).Value = "  -3.08%  "
